# Add a new worksheet "cloud_run" after the existing "compute_instance"
# sheet, cloned from it, with the project id in D3 changed to the
# cloud-run specific value and the active cell moved to D4.

$wb = $excel.ActiveWorkbook
$src = $wb.Worksheets.Item("compute_instance")

# Insert the new sheet right after the last existing sheet so it lands
# at the end of the tab strip (matches sheetId="3" after sheetId="2").
$dst = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$dst.Name = "cloud_run"

# Copy the full used range (values, formulas, number formats and styles).
$src.Range("A1:H19").Copy($dst.Range("A1"))

# Column widths match the source sheet's custom widths.
$dst.Columns("A").ColumnWidth = 22
$dst.Columns("B").ColumnWidth = 16.81640625
$dst.Columns("C").ColumnWidth = 13.453125
$dst.Columns("D").ColumnWidth = 31
$dst.Columns("F").ColumnWidth = 25.453125
$dst.Columns("G").ColumnWidth = 28.26953125
$dst.Columns("H").ColumnWidth = 18.453125

# Row heights for wrapped multi-line cells.
$dst.Rows("2").RowHeight = 29
$dst.Rows("7").RowHeight = 116
$dst.Rows("8").RowHeight = 232
$dst.Rows("10").RowHeight = 29
$dst.Rows("11").RowHeight = 29
$dst.Rows("12").RowHeight = 29
$dst.Rows("13").RowHeight = 145
$dst.Rows("14").RowHeight = 43.5
$dst.Rows("15").RowHeight = 58
$dst.Rows("16").RowHeight = 29
$dst.Rows("17").RowHeight = 29
$dst.Rows("18").RowHeight = 72.5
$dst.Rows("19").RowHeight = 58

# This sheet is for the "cloud_run" variant, so give it its own sample
# project/instance-name value instead of the compute_instance one.
$dst.Range("D3").Value = "cloud-run-01"

# Select D4 and make this new sheet the active tab, like the author left it.
$dst.Range("D4").Select()
$dst.Activate()
